$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# A new weekly price record needs to be inserted right after the existing
# header/top rows (row 4), pushing all the subsequent records down by one
# row - mirroring how the source data feed prepends its newest entry.
$ws.Rows.Item(4).Insert()

$ws.Cells.Item(4, 1).Value2 = 8
$ws.Cells.Item(4, 2).Value2 = "Terminal La Palmera de La Serena"
$ws.Cells.Item(4, 3).Value2 = "Coquimbo"
$ws.Cells.Item(4, 4).Value2 = 45043
$ws.Cells.Item(4, 5).Value2 = 4
$ws.Cells.Item(4, 6).Value2 = 100112026
$ws.Cells.Item(4, 7).Value2 = "Haba"
$ws.Cells.Item(4, 8).Value2 = "Sin especificar"
$ws.Cells.Item(4, 9).Value2 = "Primera"
$ws.Cells.Item(4, 10).Value2 = 320
$ws.Cells.Item(4, 11).Value2 = 14000
$ws.Cells.Item(4, 12).Value2 = 15000
$ws.Cells.Item(4, 13).Value2 = 14500
$ws.Cells.Item(4, 14).Value2 = "`$/saco 25 kilos"
$ws.Cells.Item(4, 15).Value2 = "Provincia del Elquí"
$ws.Cells.Item(4, 16).Value2 = 580
$ws.Cells.Item(4, 17).Value2 = 25
$ws.Cells.Item(4, 18).Value2 = "Hortaliza"
